$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (Insurance Exp / Permit Exp) to include date format hint
$ws.Range("G1").Value = "Insurance Exp (YYYY-MM-DD)"
$ws.Range("H1").Value = "Permit Exp (YYYY-MM-DD)"

# Update number format for columns G and H from text to date (built-in numFmtId 14 = mm-dd-yy)
$ws.Range("H1").NumberFormat = "mm-dd-yy"
$ws.Range("H1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update column widths (G, H widened; new I column added)
$ws.Columns("G").ColumnWidth = 25
$ws.Columns("H").ColumnWidth = 22.333333333333332
$ws.Columns("I").ColumnWidth = 8

# Update the active selection shown on the sheet
$ws.Range("H1:H1048576").Select()
